$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-5
# from 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185)
$newDate = (Get-Date -Year 2023 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
